# "Table sent to HTML"
# The "Hours" sheet's column B used to hold the resource's "Business Unit"
# (text, looked up by the hidden "Totals (2)" sheet through a VLOOKUP).
# It is repurposed to hold the resource's numeric "Rate" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hours")

# Header: "Business Unit" -> "Rate"
$ws.Range("B1").Value = "Rate"

# Resource rate values (row 2..20), replacing the old Business Unit text.
$rates = @{
    2  = 96
    3  = 58
    4  = 48
    5  = 58
    6  = 58
    7  = 48
    8  = 58
    9  = 48
    10 = 58
    11 = 48
    12 = 58
    13 = 58
    14 = 48
    15 = 58
    16 = 48
    17 = 58
    18 = 58
    19 = 48
    20 = 58
}

foreach ($row in $rates.Keys) {
    $ws.Cells.Item($row, 2).Value = $rates[$row]
}

# The downstream "Totals (2)" sheet's Business Unit column (Table22) looks
# the value up from Hours!A:B via VLOOKUP, so it recalculates automatically
# once Hours!B holds numbers instead of text.

# Restore the active selection left on the sheet (B18:B20, active cell B18).
$ws.Activate() | Out-Null
$ws.Range("B18:B20").Select() | Out-Null
